$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-08 18:15:20", 0.0004),
    @("2023-12-08 18:15:41", 0.0016),
    @("2023-12-08 18:15:59", 0.0008),
    @("2023-12-08 18:16:08", 0.0006000000000000001),
    @("2023-12-08 18:16:18", 0.0004)
)

$startRow = 120
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
